# Update the dSF (column F) values for several rows as part of a
# "repull data, push all data, mean calculation" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = 2
    5  = 0
    10 = -3
    12 = -1
    20 = 5
    21 = -11
    24 = -1
    27 = 4
    29 = -5
    30 = -3
    31 = -11
    34 = -7
    35 = -4
    38 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
